$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Justifications (if any)" column is being added as column H.
# Give H1:H8 the same look as the neighbouring column G (header style for
# row 1, bordered data style for rows 2-8) before filling in the values.
$ws.Range("G1:G8").Copy()
$ws.Range("H1:H8").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "Justifications (if any)"
for ($r = 2; $r -le 8; $r++) {
    $ws.Range("H$r").Value = "-"
}

# Widen column H so the longer header text fits comfortably.
$ws.Columns.Item(8).ColumnWidth = 24.42

# Leave the selection where the author ended up after adding the column.
$ws.Range("H9").Select()
